$wb = $excel.ActiveWorkbook

# Rename the "Include from Binary Data Enco" sheet to "Include #0"
$includeSheet = $wb.Worksheets.Item("Include from Binary Data Enco")
$includeSheet.Name = "Include #0"

# Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 (shifts existing rows 11-14 down to 12-15)
$ws.Rows.Item(11).Insert()

# Copy style from the row above (row 10) into the new row 11 cells
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set new Jurisdiction row content
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
